# Update the handoff/handback datetimes for the last row (6966655e-...) on the
# zh-cn and de-de report sheets, as produced by a fresh "Generate Report for
# Handback" run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D6").Value = "2016-03-01 09:08:35"
$zhcn.Range("G6").Value = "2016-03-01 09:09:20"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D6").Value = "2016-03-01 09:08:46"
$dede.Range("G6").Value = "2016-03-01 09:09:40"
